$d = $word.ActiveDocument

# --- Paragraph 7: "...a n by n table..." -> wrap "n by n" with gramStart/gramEnd proofErr markers ---
$p7xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w14:paraId="7A463400" w14:textId="4234715C" w:rsidR="00207437" w:rsidRDefault="00E113BD" w:rsidP="00207437" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:pPr><w:ind w:firstLine="360"/><w:jc w:val="both"/></w:pPr><w:r><w:t>Th</w:t></w:r><w:r w:rsidR="00207437"><w:t xml:space="preserve">e updated greedy </w:t></w:r><w:r><w:t xml:space="preserve">algorithm works similarly to the </w:t></w:r><w:r w:rsidR="00207437"><w:t xml:space="preserve">basic </w:t></w:r><w:r><w:t xml:space="preserve">greedy algorithm described </w:t></w:r><w:r w:rsidR="00207437"><w:t>above but</w:t></w:r><w:r><w:t xml:space="preserve"> includes backtracking.  It keeps track of which routes have been tried from a node in a </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>n by n</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> table (n is the number of cities) and tries new paths in order of shortest distance.  If there is not a valid </w:t></w:r><w:r w:rsidR="00207437"><w:t>complete tour</w:t></w:r><w:r><w:t xml:space="preserve">, this algorithm will try every possible </w:t></w:r><w:r w:rsidR="00207437"><w:t>route.</w:t></w:r></w:p>'

$p7 = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text -like "*n by n table*") {
        $p7 = $cand
        break
    }
}
if ($p7 -eq $null) {
    throw "Could not locate the 'n by n table' paragraph."
}
$null = $p7.Range.InsertXML($p7xml)

# --- Paragraph 8: "...O(n!*n2)..." appears twice -> wrap "n!*" with gramStart/gramEnd proofErr markers ---
$p8xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w14:paraId="623B9D8B" w14:textId="5953CC0C" w:rsidR="00E62D13" w:rsidRDefault="00E113BD" w:rsidP="00207437" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:pPr><w:ind w:firstLine="360"/><w:jc w:val="both"/></w:pPr><w:r><w:t>It has space complexity of O(n</w:t></w:r><w:r><w:rPr><w:vertAlign w:val="superscript"/></w:rPr><w:t>2</w:t></w:r><w:r><w:t>), best-case time complexity of O(n</w:t></w:r><w:r><w:rPr><w:vertAlign w:val="superscript"/></w:rPr><w:t>2</w:t></w:r><w:r><w:t>), and worst-case time complexity of O(</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>n!*</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>n</w:t></w:r><w:r><w:rPr><w:vertAlign w:val="superscript"/></w:rPr><w:t>2</w:t></w:r><w:r><w:t>).  It is possible to reduce the worst-case complexity to O(</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>n!*</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>n) by first</w:t></w:r><w:r w:rsidR="00207437"><w:t xml:space="preserve"> checking for a valid path without taking the time to compare distances, but the graphs we’re working with are sufficiently dense that we never found one without at least one complete tour.  We also felt it was best to leave that feature out because for best-case runs, which are common, it doubles the time taken.</w:t></w:r></w:p>'

$p8 = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text -like "*worst-case time complexity*") {
        $p8 = $cand
        break
    }
}
if ($p8 -eq $null) {
    throw "Could not locate the 'worst-case time complexity' paragraph."
}
$null = $p8.Range.InsertXML($p8xml)

# --- Last paragraph (empty, under "Discussion of Table"): add first-line indent + new body text ---
$p16xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w14:paraId="10CD6ED1" w14:textId="77777777" w:rsidR="00E62D13" w:rsidRPr="00E62D13" w:rsidRDefault="00E62D13" w:rsidP="00207437" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:pPr><w:ind w:firstLine="360"/><w:jc w:val="both"/></w:pPr><w:r><w:t xml:space="preserve">We took the average of 5 trials for city sizes of 15, 30, 60, 100, 200, 500, 100, 1500, and 2000 for each algorithm (Random, Greedy, Branch and Bound, and Local Search). Of all the algorithms, Greedy was by </w:t></w:r><w:r><w:lastRenderedPageBreak/><w:t>far the fastest</w:t></w:r><w:r><w:t xml:space="preserve"> across all city sizes. Greedy finished</w:t></w:r><w:r><w:t xml:space="preserve"> with an average run time of 10 seconds for 2000 cities. The only other algorithm that could finish a solution before the 10-minute limit was the Local Search algorithm with an average run time of 424 seconds. Nonetheless, the Local Search algorithm had the best average cost of tour for each tried city size. The Local Search algorithm</w:t></w:r><w:r><w:t xml:space="preserve"> on average had a cost of tour score 6% lower than Greedy</w:t></w:r><w:r><w:t>.</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">Therefore, Local Search is the more effective algorithm at finding the best minimum cost solution while Greedy finds good, </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>low cost</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> solutions quickly. </w:t></w:r></w:p>'

$p16 = $d.Paragraphs.Item($d.Paragraphs.Count)
if ($p16.Range.Text.Trim() -ne "") {
    throw "Expected the final paragraph to be empty before inserting the results discussion text."
}
$null = $p16.Range.InsertXML($p16xml)

Write-Host "Done."
